# Adds an "audio_name" (mp3) column entry to several rows of the "safety"
# worksheet, mirroring the existing pattern already used on the "departure"
# worksheet (where column D stores the audio file name for each row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("safety")

# New column D needs the same custom-width treatment as the sheet's other
# data columns (B and C already have explicit widths).
$ws.Columns.Item(4).ColumnWidth = 45.86

# xlPasteFormats: paste only the formatting (font/fill/alignment/etc.) of a
# template cell onto the freshly populated audio cell, so the new cell ends
# up re-using the same style as its neighbours instead of a brand-new one.
$xlPasteFormats = -4122

# sa_001 (row 2) - top level entry, style matches column C on the same row
$ws.Range("D2").Value = "right_and_your_safety.mp3"
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial($xlPasteFormats)

# sa_001_1 (row 3)
$ws.Range("D3").Value = "violence_on_women.mp3"
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial($xlPasteFormats)

# sa_001_2 (row 4) - uses the alternate "highlighted" style, like C7
$ws.Range("D4").Value = "harassment.mp3"
$ws.Range("C7").Copy()
$ws.Range("D4").PasteSpecial($xlPasteFormats)

# sa_001_3 (row 5)
$ws.Range("D5").Value = "your_health.mp3"
$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial($xlPasteFormats)

# sa_001_4 (row 6) - no image column populated, style matches B11's style
$ws.Range("D6").Value = "understand_exploitation_and_human_traficking.mp3"
$ws.Range("B11").Copy()
$ws.Range("D6").PasteSpecial($xlPasteFormats)

# sa_002 (row 8) - top level entry
$ws.Range("D8").Value = "my_body_my_choice_safety_planning_tips.mp3"
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial($xlPasteFormats)

$ws.Range("A1").Select()
